$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price values using "." as both thousands and decimal
# separators (e.g. "27.874.38"), which must stay text, not be reparsed
# as numbers by Excel (which would also strip significant trailing
# zeroes, e.g. "16.20" -> 16.2). Force text format before assignment.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.874.38"
$ws.Range("E2").Value = "  +2.71%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.665.72"
$ws.Range("E3").Value = "  -0.74%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "214.48"
$ws.Range("E5").Value = "  +0.11%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.515"
$ws.Range("E6").Value = "  -0.69%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "23.49"
$ws.Range("E8").Value = "  +3.06%  "
$ws.Range("E9").Value = "  +0.36%  "
$ws.Range("E10").Value = "  -0.32%  "
$ws.Range("E11").Value = "  -1.32%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.898.00"
$ws.Range("E12").Value = "  -0.87%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.658.17"
$ws.Range("E13").Value = "  -1.10%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.14"
$ws.Range("E14").Value = "  -1.54%  "
$ws.Range("E15").Value = "  -0.40%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.83"
$ws.Range("E16").Value = "  -1.07%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "249.35"
$ws.Range("E17").Value = "  +6.13%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "27.826.01"
$ws.Range("E18").Value = "  +2.67%  "
$ws.Range("E19").Value = "  -1.17%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.54"
$ws.Range("E20").Value = "  -4.43%  "
$ws.Range("E21").Value = "  +0.21%  "
$ws.Range("E22").Value = "  -1.56%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.35"
$ws.Range("E23").Value = "  -1.73%  "
$ws.Range("E24").Value = "  -1.62%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "146.77"
$ws.Range("E25").Value = "  -1.05%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.21"
$ws.Range("E26").Value = "  -2.87%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.20"
$ws.Range("E27").Value = "  -0.85%  "
$ws.Range("E28").Value = "  +0.01%  "
$ws.Range("E29").Value = "  -0.87%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.23"
$ws.Range("E30").Value = "  +5.70%  "
$ws.Range("E31").Value = "  -0.52%  "
$ws.Range("E32").Value = "  -0.63%  "
$ws.Range("E33").Value = "  -3.09%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.414.35"
$ws.Range("E34").Value = "  -8.15%  "
$ws.Range("E35").Value = "  -5.70%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.39"
$ws.Range("E36").Value = "  +0.02%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.931"
$ws.Range("E37").Value = "  -0.76%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.578"
$ws.Range("E38").Value = "  -4.45%  "
$ws.Range("E39").Value = "  -1.81%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.04"
$ws.Range("E40").Value = "  -2.53%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "69.20"
$ws.Range("E41").Value = "  -0.26%  "
$ws.Range("E42").Value = "  +0.02%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.43"
$ws.Range("E43").Value = "  -6.16%  "
$ws.Range("E44").Value = "  -1.16%  "
$ws.Range("E45").Value = "  +1.63%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.806.72"
$ws.Range("E46").Value = "  -0.87%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.71"
$ws.Range("E47").Value = "  +4.78%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "88.36"
$ws.Range("E48").Value = "  -1.72%  "
$ws.Range("E49").Value = "  -4.14%  "
$ws.Range("E50").Value = "  -2.55%  "
$ws.Range("B51").NumberFormat = "@"
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").NumberFormat = "@"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.77"
$ws.Range("E51").Value = "  -5.13%  "
